$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values to match the re-pulled/recalculated data
$ws.Range("F2").Value = -1
$ws.Range("F6").Value = -4
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = 4
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = 1
